$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.622.58"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "1.926.96"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "'326.36"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("D7").Value = "'0.4819"
$ws.Range("E7").Value = "  -0.54%  "
$ws.Range("D8").Value = "'0.4058"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "'0.08199"
$ws.Range("E9").Value = "  +0.28%  "
$ws.Range("D10").Value = "'1.009"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").Value = "'23.80"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'6.077"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'7.294"
$ws.Range("E13").Value = "  +1.43%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.897.10"
$ws.Range("E14").Value = "  -0.65%  "
$ws.Range("D15").Value = "'91.52"
$ws.Range("E15").Value = "  +0.40%  "
$ws.Range("D16").Value = "'0.06854"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "'17.58"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "'1.010"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "29.600.84"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'5.657"
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("D23").Value = "'11.96"
$ws.Range("E23").Value = "  +1.81%  "
$ws.Range("D24").Value = "'2.185"
$ws.Range("E24").Value = "  +0.35%  "
$ws.Range("D25").Value = "2.075.15"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").Value = "'156.30"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'6.360"
$ws.Range("E27").Value = "  -2.50%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'20.05"
$ws.Range("E28").Value = "  -0.37%  "
$ws.Range("D29").Value = "'2.089"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "'120.81"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'1.005"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").Value = "'0.09602"
$ws.Range("E32").Value = "  +0.49%  "
$ws.Range("D33").Value = "'5.619"
$ws.Range("E33").Value = "  +2.04%  "
$ws.Range("D34").Value = "'3.557"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'1.391"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").Value = "'0.06535"
$ws.Range("E36").Value = "  +6.63%  "
$ws.Range("D37").Value = "'0.02280"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "'1.220"
$ws.Range("E38").Value = "  +2.71%  "
$ws.Range("D39").Value = "'0.5942"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "'10.73"
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'7.869"
$ws.Range("E41").Value = "  -1.61%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1844"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'2.495"
$ws.Range("E43").Value = "  +3.37%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.245"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").Value = "'0.07546"
$ws.Range("E45").Value = "  -1.07%  "
$ws.Range("D46").Value = "'12.34"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5555"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.956"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'118.66"
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.430"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'71.93"
$ws.Range("E51").Value = "  -1.00%  "
